# Scheduled market-data refresh: update cached Leve profit values across sheets
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 100000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 100000
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H135").Value = 521.76
$ws.Range("J135").Value = 727.3333
$ws.Range("L135").Value = 6545.9997
$ws.Range("N135").Value = -11615.9997
# Leve with no remaining HQ price data point: clear stale profit figure
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9102.9
$ws.Range("J2").Value = 19104.2
$ws.Range("L2").Value = 19104.2
$ws.Range("N2").Value = -19330.2
$ws.Range("H4").Value = 174.90909
$ws.Range("J4").Value = 93.5
$ws.Range("L4").Value = 93.5
$ws.Range("N4").Value = -325.5
$ws.Range("H45").Value = 3280.7144
$ws.Range("I45").Value = 2910.4666
$ws.Range("K45").Value = 2910.4666
$ws.Range("M45").Value = -2533.4666
$ws.Range("H61").Value = 3335202
$ws.Range("I61").Value = 3705446.8
$ws.Range("K61").Value = 3705446.8
$ws.Range("M61").Value = -3705234.8
$ws.Range("H63").Value = 7522.647
$ws.Range("I63").Value = 3898
$ws.Range("K63").Value = 3898
$ws.Range("M63").Value = -3212
$ws.Range("H66").Value = 7522.647
$ws.Range("I66").Value = 3898
$ws.Range("K66").Value = 19490
$ws.Range("M66").Value = -16058
$ws.Range("H74").Value = 2971.8276
$ws.Range("I74").Value = 1634.25
$ws.Range("K74").Value = 1634.25
$ws.Range("M74").Value = -760.25
$ws.Range("H77").Value = 2971.8276
$ws.Range("I77").Value = 1634.25
$ws.Range("K77").Value = 8171.25
$ws.Range("M77").Value = -3803.25
$ws.Range("H102").Value = 150448
$ws.Range("I102").Value = 150448
$ws.Range("K102").Value = 150448
$ws.Range("M102").Value = -148826
$ws.Range("H116").Value = 9102.9
$ws.Range("J116").Value = 19104.2
$ws.Range("L116").Value = 19104.2
$ws.Range("N116").Value = -23692.2
$ws.Range("H136").Value = 3335202
$ws.Range("I136").Value = 3705446.8
$ws.Range("K136").Value = 11116340.4
$ws.Range("M136").Value = -11113790.4

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9102.9
$ws.Range("J3").Value = 19104.2
$ws.Range("L3").Value = 19104.2
$ws.Range("N3").Value = -19332.2
$ws.Range("H20").Value = 972.0714
$ws.Range("I20").Value = 966.9048
$ws.Range("J20").Value = 987.5714
$ws.Range("K20").Value = 966.9048
$ws.Range("L20").Value = 987.5714
$ws.Range("M20").Value = -719.9048
$ws.Range("N20").Value = -1481.5714
$ws.Range("H22").Value = 311.35715
$ws.Range("I22").Value = 354.75
$ws.Range("J22").Value = 51
$ws.Range("K22").Value = 354.75
$ws.Range("L22").Value = 51
$ws.Range("M22").Value = -181.75
$ws.Range("N22").Value = -397
$ws.Range("H105").Value = 4833.3335
$ws.Range("I105").Value = 4833.3335
$ws.Range("K105").Value = 4833.3335
$ws.Range("M105").Value = -3086.3335
$ws.Range("H107").Value = 2624
$ws.Range("I107").Value = 2016.72
$ws.Range("K107").Value = 2016.72
$ws.Range("M107").Value = -96.72000000000003

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 738.3333
$ws.Range("I16").Value = 486
$ws.Range("K16").Value = 486
$ws.Range("M16").Value = -199
$ws.Range("H31").Value = 23636.041
$ws.Range("I31").Value = 14152.857
$ws.Range("K31").Value = 14152.857
$ws.Range("M31").Value = -13857.857
$ws.Range("H34").Value = 23636.041
$ws.Range("I34").Value = 14152.857
$ws.Range("K34").Value = 14152.857
$ws.Range("M34").Value = -13950.857
$ws.Range("H62").Value = 8586.286
$ws.Range("J62").Value = 9184.333000000001
$ws.Range("L62").Value = 9184.333000000001
$ws.Range("N62").Value = -10432.333
$ws.Range("H65").Value = 8586.286
$ws.Range("J65").Value = 9184.333000000001
$ws.Range("L65").Value = 45921.665
$ws.Range("N65").Value = -52161.665
$ws.Range("H113").Value = 738.3333
$ws.Range("I113").Value = 486
$ws.Range("K113").Value = 486
$ws.Range("M113").Value = 1684
$ws.Range("H117").Value = 30000
$ws.Range("I117").Value = 30000
$ws.Range("K117").Value = 30000
$ws.Range("M117").Value = -25411
$ws.Range("H132").Value = 2696.2334
$ws.Range("I132").Value = 1878.4736
$ws.Range("K132").Value = 5635.4208
$ws.Range("M132").Value = -3105.4208

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 40423890
$ws.Range("I4").Value = 46250230
$ws.Range("J4").Value = 21002750
$ws.Range("K4").Value = 138750690
$ws.Range("L4").Value = 63008250
$ws.Range("M4").Value = -138750578
$ws.Range("N4").Value = -63008474
$ws.Range("H32").Value = 650794500
$ws.Range("J32").Value = 650794500
$ws.Range("L32").Value = 1952383500
$ws.Range("N32").Value = -1952384066
$ws.Range("H92").Value = 685.2857
$ws.Range("J92").Value = 739
$ws.Range("L92").Value = 2217
$ws.Range("N92").Value = -4713
$ws.Range("H107").Value = 899.1667
$ws.Range("J107").Value = 918.3333
$ws.Range("L107").Value = 2754.9999
$ws.Range("N107").Value = -6594.9999
$ws.Range("H133").Value = 4276.6665
$ws.Range("I133").Value = 4276.6665
$ws.Range("K133").Value = 12829.9995
$ws.Range("M133").Value = -7769.999500000002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1891.25
$ws.Range("I97").Value = 1963.7693
$ws.Range("J97").Value = 1577
$ws.Range("K97").Value = 1963.7693
$ws.Range("L97").Value = 1577
$ws.Range("M97").Value = -1467.7693
$ws.Range("N97").Value = -2569
$ws.Range("H102").Value = 2804.1052
$ws.Range("I102").Value = 1819.5385
$ws.Range("J102").Value = 4937.3335
$ws.Range("K102").Value = 1819.5385
$ws.Range("L102").Value = 4937.3335
$ws.Range("M102").Value = -197.5385000000001
$ws.Range("N102").Value = -8181.3335
$ws.Range("H122").Value = 59503.777
$ws.Range("I122").Value = 62710.47
$ws.Range("K122").Value = 188131.41
$ws.Range("M122").Value = -185681.41

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3953.5217
$ws.Range("I7").Value = 3808.6875
$ws.Range("K7").Value = 3808.6875
$ws.Range("M7").Value = -3696.6875
$ws.Range("H40").Value = 2220.0344
$ws.Range("I40").Value = 2233.7307
$ws.Range("K40").Value = 2233.7307
$ws.Range("M40").Value = -2097.7307
$ws.Range("H61").Value = 2067.7144
$ws.Range("I61").Value = 1754.2354
$ws.Range("K61").Value = 1754.2354
$ws.Range("M61").Value = -1552.2354
$ws.Range("H113").Value = 2067.7144
$ws.Range("I113").Value = 1754.2354
$ws.Range("K113").Value = 1754.2354
$ws.Range("M113").Value = 415.7646
$ws.Range("H114").Value = 100318.6
$ws.Range("J114").Value = 100318.6
$ws.Range("L114").Value = 100318.6
$ws.Range("N114").Value = -108996.6
$ws.Range("H126").Value = 3953.5217
$ws.Range("I126").Value = 3808.6875
$ws.Range("K126").Value = 11426.0625
$ws.Range("M126").Value = -8956.0625
$ws.Range("H130").Value = 85653.28999999999
$ws.Range("J130").Value = 85653.28999999999
$ws.Range("L130").Value = 85653.28999999999
$ws.Range("N130").Value = -95693.28999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1047.1305
$ws.Range("I107").Value = 640.63635
$ws.Range("K107").Value = 1921.90905
$ws.Range("M107").Value = -1.909049999999979
$ws.Range("H122").Value = 1718.7931
$ws.Range("I122").Value = 1392.9048
$ws.Range("K122").Value = 4178.7144
$ws.Range("M122").Value = -1728.7144
$ws.Range("H126").Value = 1894.1818
$ws.Range("I126").Value = 1562.6316
$ws.Range("J126").Value = 3994
$ws.Range("K126").Value = 4687.8948
$ws.Range("L126").Value = 11982
$ws.Range("M126").Value = -2217.8948
$ws.Range("N126").Value = -16922
$ws.Range("H132").Value = 2489695.8
$ws.Range("I132").Value = 3032232.8
$ws.Range("K132").Value = 9096698.399999999
$ws.Range("M132").Value = -9094168.399999999
$ws.Range("H136").Value = 9479.799999999999
$ws.Range("I136").Value = 9374.791999999999
$ws.Range("J136").Value = 12000
$ws.Range("K136").Value = 28124.376
$ws.Range("L136").Value = 36000
$ws.Range("M136").Value = -25574.376
$ws.Range("N136").Value = -41100
